$wb = $excel.ActiveWorkbook

# --- Settings sheet: add new configuration rows ---
$wsSettings = $wb.Worksheets.Item("Settings")

$wsSettings.Cells.Item(6, 1).Value = "SheetName"
$wsSettings.Cells.Item(6, 2).Value = "Sheet1"

$wsSettings.Cells.Item(7, 1).Value = "Orchestrator_FolderName"
$wsSettings.Cells.Item(7, 2).Value = "Robot1"

$wsSettings.Cells.Item(8, 1).Value = "RPAChallenge_URL"
$wsSettings.Cells.Item(8, 2).Value = "RPAChallenge_URL"

$wsSettings.Cells.Item(9, 1).Value = "RPAChallenge_Path"
$wsSettings.Cells.Item(9, 2).Value = "RPAChallenge_Path"

$newRowRange = $wsSettings.Range("A9:B9")
$newRowRange.Font.Name = "Segoe UI"
$newRowRange.Font.Color = 5590598

# --- Assets sheet: move selection to B11 without leaving it as the active sheet ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("B11").Select()

# --- Restore Settings as the active/tab-selected sheet, with its new selection ---
$wsSettings.Activate()
$wsSettings.Range("A18").Select()
